# Update "想去人数" (attendance interest count) figures on the "展览" and
# "全部类型" worksheets, reflecting the freshly re-generated data snapshot.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> cell address -> new value
$updates = @{
    "展览" = @{
        "F2"  = 123
        "F3"  = 2141
        "F5"  = 11182
        "F8"  = 310
        "F10" = 11090
        "F13" = 43
        "F14" = 1726
        "F15" = 5561
        "F16" = 92
        "F17" = 3438
    }
    "全部类型" = @{
        "F2"  = 123
        "F3"  = 2141
        "F7"  = 11182
        "F10" = 310
        "F12" = 11090
        "F15" = 43
        "F16" = 1726
        "F17" = 5561
        "F18" = 92
        "F19" = 3438
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($addr in $cellUpdates.Keys) {
        $ws.Range($addr).Value = $cellUpdates[$addr]
    }
}
